$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 6 was previously blank (A6 had a style only, B6 had nothing).
# Restore the filter values: A6 = "heidelberg-03", B6 = "NC_011083.1-602044"
# Set B6 first so the shared-string table records "NC_011083.1-602044"
# before "heidelberg-03", matching the original authoring order.
$ws.Range("B6").Value = "NC_011083.1-602044"
$ws.Range("A6").Value = "heidelberg-03"
